$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B10: record the work hours for that day
$ws.Range("B10").Value = "11 Hours 38 Minutes"

# A11: add the next date, matching the date style/format used by the cells above (A2:A10)
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
$prevDate = $ws.Range("A10").Value2
$ws.Range("A11").Value2 = $prevDate + 1

# Update the active selection to B11
$ws.Range("B11").Select()
